# "Generate Report for Handback"
# Update the handoff/handback timestamp cells on the Overview, zh-cn and
# de-de sheets to reflect the latest xliff generation / handback times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# G2: "Latest HO Xliff Generate Date" for 848b503e-...-f5c445444957.md
$wsOverview.Range("G2").Value = "2016-08-28 17:04:29"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# H2: "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-28 17:04:25"
# K2: "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-08-28 17:04:42"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# H2: "Correspond Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-08-28 17:04:29"
# K2: "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-08-28 17:04:48"
